# update scripts wuth new tpm
# Refresh the NATMI ligand-receptor statistics in Fn1-Itgb7.xlsx with
# values recomputed from the new TPM expression matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("M2").Value = 0.8194946666666666
$ws.Range("N2").Value = 2.458484
$ws.Range("O2").Value = 0.1466535424263973
$ws.Range("P2").Value = 0.1466535424263973
$ws.Range("Q2").Value = 5.113735225423999
$ws.Range("R2").Value = 46.02361702881599
$ws.Range("S2").Value = 0.002540377425761823
$ws.Range("T2").Value = 0.002540377425761824

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.6452529427684778
$ws.Range("P3").Value = 0.6452529427684778
$ws.Range("Q3").Value = 22.49964541019999
$ws.Range("R3").Value = 202.4968086917999
$ws.Range("S3").Value = 0.01117726842867163
$ws.Range("T3").Value = 0.01117726842867164

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("N4").Value = 3.488457
$ws.Range("O4").Value = 0.2080935148051249
$ws.Range("P4").Value = 0.2080935148051249
$ws.Range("Q4").Value = 7.256116144451999
$ws.Range("R4").Value = 65.30504530006799
$ws.Range("S4").Value = 0.003604659380960304
$ws.Range("T4").Value = 0.003604659380960305

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("M5").Value = 0.8194946666666666
$ws.Range("N5").Value = 2.458484
$ws.Range("O5").Value = 0.1466535424263973
$ws.Range("P5").Value = 0.1466535424263973
$ws.Range("Q5").Value = 283.1899684687453
$ws.Range("R5").Value = 2548.709716218708
$ws.Range("S5").Value = 0.1406817856981545
$ws.Range("T5").Value = 0.1406817856981546

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.6452529427684778
$ws.Range("P6").Value = 0.6452529427684778
$ws.Range("S6").Value = 0.6189781352279098
$ws.Range("T6").Value = 0.6189781352279099

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("N7").Value = 3.488457
$ws.Range("O7").Value = 0.2080935148051249
$ws.Range("P7").Value = 0.2080935148051249
$ws.Range("Q7").Value = 401.831383826201
$ws.Range("S7").Value = 0.1996199121455446
$ws.Range("T7").Value = 0.1996199121455446

# Row 8 (MuSCs -> ECs)
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("M8").Value = 0.8194946666666666
$ws.Range("N8").Value = 2.458484
$ws.Range("O8").Value = 0.1466535424263973
$ws.Range("P8").Value = 0.1466535424263973
$ws.Range("Q8").Value = 6.907306384060443
$ws.Range("R8").Value = 62.16575745654399
$ws.Range("S8").Value = 0.003431379302480951
$ws.Range("T8").Value = 0.003431379302480953

# Row 9 (MuSCs -> FAPs)
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.6452529427684778
$ws.Range("P9").Value = 0.6452529427684778
$ws.Range("R9").Value = 273.5197341611999
$ws.Range("S9").Value = 0.01509753911189633
$ws.Range("T9").Value = 0.01509753911189633

# Row 10 (MuSCs -> MuSCs)
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("N10").Value = 3.488457
$ws.Range("O10").Value = 0.2080935148051249
$ws.Range("P10").Value = 0.2080935148051249
$ws.Range("S10").Value = 0.004868943278619992
$ws.Range("T10").Value = 0.004868943278619994
